# Logic tree input file updated
# Insert two new "Possible_Problem" placeholder rows into the decision
# tree: one right above the "Normal/Cold/Hot(Red)/Fluctuates" answer
# block, and one right above the "Yes/No" (fan loudness) answer block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture text we will reuse, from the ORIGINAL (pre-insert) layout ---
$txtQ1       = $ws.Range("A4").Value2   # "Problem:What is shown on the engine temp gauge..."
$txtPossible = $ws.Range("B8").Value2   # "Possible_Problem"
$txtLeakRad  = $ws.Range("C7").Value2   # "Possible_Problem:16% Leaking Radiator..."
$txtQ4       = $ws.Range("A19").Value2  # "Problem:Does the fan become louder..."
$txtDup      = $ws.Range("C12").Value2  # "Possible_Problem:33% Plugged heater core..."

# --- insert new row 4 (above the Normal/Cold/Hot/Fluctuates answers) ---
$ws.Rows(4).Insert()
$ws.Range("A4").Value = $txtQ1
$ws.Range("B4").Value = $txtPossible
$ws.Range("C4").Value = $txtLeakRad
$ws.Rows(4).RowHeight = 409.6
$ws.Range("C4").WrapText = $true

# --- insert new row 20 (above the Yes/No fan-loudness answers) ---
$ws.Rows(20).Insert()
$ws.Range("A20").Value = $txtQ4
$ws.Range("B20").Value = $txtPossible
$ws.Range("C20").Value = $txtDup
$ws.Rows(20).RowHeight = 403.2
$ws.Range("C20").WrapText = $true

# --- refresh the view to focus the new last data rows, matching the saved file ---
$ws.Range("A22").Select()
$ws.Application.ActiveWindow.ScrollRow = 22
